$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 10. This shifts the existing data rows
# 10..83 down to 11..84 (which is exactly what the target diff shows:
# every existing record moves down by one row, and a brand new record
# is introduced at row 10).
$ws.Rows(10).Insert()

# Populate the newly inserted row 10 with the new record's data.
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C10").Value = "Metropolitana"
$ws.Range("D10").Value = 45149
$ws.Range("E10").Value = 13
$ws.Range("F10").Value = 100112010
$ws.Range("G10").Value = "Achicoria"
$ws.Range("H10").Value = "Sin especificar"
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 90
$ws.Range("K10").Value = 7000
$ws.Range("L10").Value = 7000
$ws.Range("M10").Value = 7000
$ws.Range("N10").Value = '$/caja 16 unidades'
$ws.Range("O10").Value = "Provincia de Quillota"
$ws.Range("P10").Value = 438
$ws.Range("Q10").Value = 16
$ws.Range("R10").Value = "Hortaliza"
